$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated TPM-derived values to LR-pairs table (Ptn-Ptprs)
# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2563003333333334
$ws.Range("H2").Value = 0.7689010000000001
$ws.Range("I2").Value = 0.02986826554325775
$ws.Range("J2").Value = 0.02986826554325775
$ws.Range("M2").Value = 4.260872666666667
$ws.Range("N2").Value = 12.782618
$ws.Range("O2").Value = 0.09064705929364961
$ws.Range("P2").Value = 0.09064705929364959
$ws.Range("Q2").Value = 1.092063084757556
$ws.Range("R2").Value = 9.828567762818002
$ws.Range("S2").Value = 0.002707470437698157
$ws.Range("T2").Value = 0.002707470437698156

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2563003333333334
$ws.Range("H3").Value = 0.7689010000000001
$ws.Range("I3").Value = 0.02986826554325775
$ws.Range("J3").Value = 0.02986826554325775
$ws.Range("O3").Value = 0.4366505728284585
$ws.Range("P3").Value = 0.4366505728284584
$ws.Range("Q3").Value = 5.260512312699
$ws.Range("R3").Value = 47.344610814291
$ws.Range("S3").Value = 0.01304199525885601
$ws.Range("T3").Value = 0.01304199525885600

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2563003333333334
$ws.Range("H4").Value = 0.7689010000000001
$ws.Range("I4").Value = 0.02986826554325775
$ws.Range("J4").Value = 0.02986826554325775
$ws.Range("M4").Value = 8.931090666666666
$ws.Range("N4").Value = 26.793272
$ws.Range("O4").Value = 0.190002651698962
$ws.Range("P4").Value = 0.1900026516989619
$ws.Range("Q4").Value = 2.289041514896889
$ws.Range("R4").Value = 20.601373634072
$ws.Range("S4").Value = 0.005675049654867709
$ws.Range("T4").Value = 0.005675049654867707

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2563003333333334
$ws.Range("H5").Value = 0.7689010000000001
$ws.Range("I5").Value = 0.02986826554325775
$ws.Range("J5").Value = 0.02986826554325775
$ws.Range("M5").Value = 6.457974333333333
$ws.Range("N5").Value = 19.373923
$ws.Range("O5").Value = 0.1373888468646722
$ws.Range("P5").Value = 0.1373888468646721
$ws.Range("Q5").Value = 1.655180974291445
$ws.Range("R5").Value = 14.896628768623
$ws.Range("S5").Value = 0.004103566560836004
$ws.Range("T5").Value = 0.004103566560836001

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2563003333333334
$ws.Range("H6").Value = 0.7689010000000001
$ws.Range("I6").Value = 0.02986826554325775
$ws.Range("J6").Value = 0.02986826554325775
$ws.Range("M6").Value = 1.948535
$ws.Range("N6").Value = 5.845605
$ws.Range("O6").Value = 0.04145370713904261
$ws.Range("P6").Value = 0.0414537071390426
$ws.Range("Q6").Value = 0.4994101700116667
$ws.Range("R6").Value = 4.494691530105
$ws.Range("S6").Value = 0.001238150332581364
$ws.Range("T6").Value = 0.001238150332581364

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2563003333333334
$ws.Range("H7").Value = 0.7689010000000001
$ws.Range("I7").Value = 0.02986826554325775
$ws.Range("J7").Value = 0.02986826554325775
$ws.Range("M7").Value = 4.881814666666666
$ws.Range("N7").Value = 14.645444
$ws.Range("O7").Value = 0.1038571621752152
$ws.Range("P7").Value = 0.1038571621752152
$ws.Range("Q7").Value = 1.251210726338222
$ws.Range("R7").Value = 11.260896537044
$ws.Range("S7").Value = 0.003102033298418512
$ws.Range("T7").Value = 0.003102033298418511

# Row 8
$ws.Range("G8").Value = 6.495645000000001
$ws.Range("I8").Value = 0.7569777503270297
$ws.Range("J8").Value = 0.7569777503270296
$ws.Range("M8").Value = 4.260872666666667
$ws.Range("N8").Value = 12.782618
$ws.Range("O8").Value = 0.09064705929364961
$ws.Range("P8").Value = 0.09064705929364959
$ws.Range("Q8").Value = 27.67711623287
$ws.Range("R8").Value = 249.0940460958301
$ws.Range("S8").Value = 0.06861780701786775
$ws.Range("T8").Value = 0.06861780701786774

# Row 9
$ws.Range("G9").Value = 6.495645000000001
$ws.Range("I9").Value = 0.7569777503270297
$ws.Range("J9").Value = 0.7569777503270296
$ws.Range("O9").Value = 0.4366505728284585
$ws.Range("P9").Value = 0.4366505728284584
$ws.Range("S9").Value = 0.3305347682986954
$ws.Range("T9").Value = 0.3305347682986952

# Row 10
$ws.Range("G10").Value = 6.495645000000001
$ws.Range("I10").Value = 0.7569777503270297
$ws.Range("J10").Value = 0.7569777503270296
$ws.Range("M10").Value = 8.931090666666666
$ws.Range("N10").Value = 26.793272
$ws.Range("O10").Value = 0.190002651698962
$ws.Range("P10").Value = 0.1900026516989619
$ws.Range("Q10").Value = 58.01319443348
$ws.Range("R10").Value = 522.11874990132
$ws.Range("S10").Value = 0.1438277798392504
$ws.Range("T10").Value = 0.1438277798392504

# Row 11
$ws.Range("G11").Value = 6.495645000000001
$ws.Range("I11").Value = 0.7569777503270297
$ws.Range("J11").Value = 0.7569777503270296
$ws.Range("M11").Value = 6.457974333333333
$ws.Range("N11").Value = 19.373923
$ws.Range("O11").Value = 0.1373888468646722
$ws.Range("P11").Value = 0.1373888468646721
$ws.Range("Q11").Value = 41.948708688445
$ws.Range("R11").Value = 377.5383781960051
$ws.Range("S11").Value = 0.1040003002196443
$ws.Range("T11").Value = 0.1040003002196443

# Row 12
$ws.Range("G12").Value = 6.495645000000001
$ws.Range("I12").Value = 0.7569777503270297
$ws.Range("J12").Value = 0.7569777503270296
$ws.Range("M12").Value = 1.948535
$ws.Range("N12").Value = 5.845605
$ws.Range("O12").Value = 0.04145370713904261
$ws.Range("P12").Value = 0.0414537071390426
$ws.Range("Q12").Value = 12.656991630075
$ws.Range("R12").Value = 113.912924670675
$ws.Range("S12").Value = 0.03137953397282801
$ws.Range("T12").Value = 0.031379533972828

# Row 13
$ws.Range("G13").Value = 6.495645000000001
$ws.Range("I13").Value = 0.7569777503270297
$ws.Range("J13").Value = 0.7569777503270296
$ws.Range("M13").Value = 4.881814666666666
$ws.Range("N13").Value = 14.645444
$ws.Range("O13").Value = 0.1038571621752152
$ws.Range("P13").Value = 0.1038571621752152
$ws.Range("Q13").Value = 31.71053503046
$ws.Range("R13").Value = 285.39481527414
$ws.Range("S13").Value = 0.07861756097874388
$ws.Range("T13").Value = 0.07861756097874384

# Row 14
$ws.Range("G14").Value = 1.804372666666667
$ws.Range("H14").Value = 5.413118000000001
$ws.Range("I14").Value = 0.2102747243676212
$ws.Range("J14").Value = 0.2102747243676212
$ws.Range("M14").Value = 4.260872666666667
$ws.Range("N14").Value = 12.782618
$ws.Range("O14").Value = 0.09064705929364961
$ws.Range("P14").Value = 0.09064705929364959
$ws.Range("Q14").Value = 7.688202175880446
$ws.Range("R14").Value = 69.19381958292402
$ws.Range("S14").Value = 0.01906078540770759
$ws.Range("T14").Value = 0.01906078540770758

# Row 15
$ws.Range("G15").Value = 1.804372666666667
$ws.Range("H15").Value = 5.413118000000001
$ws.Range("I15").Value = 0.2102747243676212
$ws.Range("J15").Value = 0.2102747243676212
$ws.Range("O15").Value = 0.4366505728284585
$ws.Range("P15").Value = 0.4366505728284584
$ws.Range("Q15").Value = 37.03438269568201
$ws.Range("R15").Value = 333.309444261138
$ws.Range("S15").Value = 0.09181657884646803
$ws.Range("T15").Value = 0.09181657884646799

# Row 16
$ws.Range("G16").Value = 1.804372666666667
$ws.Range("H16").Value = 5.413118000000001
$ws.Range("I16").Value = 0.2102747243676212
$ws.Range("J16").Value = 0.2102747243676212
$ws.Range("M16").Value = 8.931090666666666
$ws.Range("N16").Value = 26.793272
$ws.Range("O16").Value = 0.190002651698962
$ws.Range("P16").Value = 0.1900026516989619
$ws.Range("Q16").Value = 16.11501588245511
$ws.Range("R16").Value = 145.035142942096
$ws.Range("S16").Value = 0.03995275521511636
$ws.Range("T16").Value = 0.03995275521511635

# Row 17
$ws.Range("G17").Value = 1.804372666666667
$ws.Range("H17").Value = 5.413118000000001
$ws.Range("I17").Value = 0.2102747243676212
$ws.Range("J17").Value = 0.2102747243676212
$ws.Range("M17").Value = 6.457974333333333
$ws.Range("N17").Value = 19.373923
$ws.Range("O17").Value = 0.1373888468646722
$ws.Range("P17").Value = 0.1373888468646721
$ws.Range("Q17").Value = 11.65259236910156
$ws.Range("R17").Value = 104.873331321914
$ws.Range("S17").Value = 0.02888940190565426
$ws.Range("T17").Value = 0.02888940190565424

# Row 18
$ws.Range("G18").Value = 1.804372666666667
$ws.Range("H18").Value = 5.413118000000001
$ws.Range("I18").Value = 0.2102747243676212
$ws.Range("J18").Value = 0.2102747243676212
$ws.Range("M18").Value = 1.948535
$ws.Range("N18").Value = 5.845605
$ws.Range("O18").Value = 0.04145370713904261
$ws.Range("P18").Value = 0.0414537071390426
$ws.Range("Q18").Value = 3.515883294043333
$ws.Range("R18").Value = 31.64294964639
$ws.Range("S18").Value = 0.008716666842678276
$ws.Range("T18").Value = 0.008716666842678274

# Row 19
$ws.Range("G19").Value = 1.804372666666667
$ws.Range("H19").Value = 5.413118000000001
$ws.Range("I19").Value = 0.2102747243676212
$ws.Range("J19").Value = 0.2102747243676212
$ws.Range("M19").Value = 4.881814666666666
$ws.Range("N19").Value = 14.645444
$ws.Range("O19").Value = 0.1038571621752152
$ws.Range("P19").Value = 0.1038571621752152
$ws.Range("Q19").Value = 8.808612948265777
$ws.Range("R19").Value = 79.27751653439199
$ws.Range("S19").Value = 0.02183853614999671
$ws.Range("T19").Value = 0.0218385361499967

# Row 20
$ws.Range("I20").Value = 0.002879259762091359
$ws.Range("J20").Value = 0.002879259762091358
$ws.Range("M20").Value = 4.260872666666667
$ws.Range("N20").Value = 12.782618
$ws.Range("O20").Value = 0.09064705929364961
$ws.Range("P20").Value = 0.09064705929364959
$ws.Range("Q20").Value = 0.1052733809753334
$ws.Range("R20").Value = 0.9474604287780002
$ws.Range("S20").Value = 0.0002609964303761148
$ws.Range("T20").Value = 0.0002609964303761148

# Row 21
$ws.Range("I21").Value = 0.002879259762091359
$ws.Range("J21").Value = 0.002879259762091358
$ws.Range("O21").Value = 0.4366505728284585
$ws.Range("P21").Value = 0.4366505728284584
$ws.Range("S21").Value = 0.001257230424439123
$ws.Range("T21").Value = 0.001257230424439122

# Row 22
$ws.Range("I22").Value = 0.002879259762091359
$ws.Range("J22").Value = 0.002879259762091358
$ws.Range("M22").Value = 8.931090666666666
$ws.Range("N22").Value = 26.793272
$ws.Range("O22").Value = 0.190002651698962
$ws.Range("P22").Value = 0.1900026516989619
$ws.Range("Q22").Value = 0.2206604571013333
$ws.Range("R22").Value = 1.985944113912
$ws.Range("S22").Value = 0.0005470669897274806
$ws.Range("T22").Value = 0.0005470669897274802

# Row 23
$ws.Range("I23").Value = 0.002879259762091359
$ws.Range("J23").Value = 0.002879259762091358
$ws.Range("M23").Value = 6.457974333333333
$ws.Range("N23").Value = 19.373923
$ws.Range("O23").Value = 0.1373888468646722
$ws.Range("P23").Value = 0.1373888468646721
$ws.Range("Q23").Value = 0.1595571718536667
$ws.Range("R23").Value = 1.436014546683
$ws.Range("S23").Value = 0.0003955781785375821
$ws.Range("T23").Value = 0.0003955781785375819

# Row 24
$ws.Range("I24").Value = 0.002879259762091359
$ws.Range("J24").Value = 0.002879259762091358
$ws.Range("M24").Value = 1.948535
$ws.Range("N24").Value = 5.845605
$ws.Range("O24").Value = 0.04145370713904261
$ws.Range("P24").Value = 0.0414537071390426
$ws.Range("Q24").Value = 0.048142454245
$ws.Range("R24").Value = 0.433282088205
$ws.Range("S24").Value = 0.0001193559909549647
$ws.Range("T24").Value = 0.0001193559909549646

# Row 25
$ws.Range("I25").Value = 0.002879259762091359
$ws.Range("J25").Value = 0.002879259762091358
$ws.Range("M25").Value = 4.881814666666666
$ws.Range("N25").Value = 14.645444
$ws.Range("O25").Value = 0.1038571621752152
$ws.Range("P25").Value = 0.1038571621752152
$ws.Range("Q25").Value = 0.1206149949693333
$ws.Range("R25").Value = 1.085534954724
$ws.Range("S25").Value = 0.0002990317480560937
$ws.Range("T25").Value = 0.0002990317480560936
